$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell = "D2"; Value = "66.098.33"}
    @{Cell = "E2"; Value = "  +0.95%  "}
    @{Cell = "D3"; Value = "3.320.75"}
    @{Cell = "E3"; Value = "  +0.50%  "}
    @{Cell = "D4"; Value = "0.999"}
    @{Cell = "E4"; Value = "  -0.06%  "}
    @{Cell = "D5"; Value = "188.19"}
    @{Cell = "E5"; Value = "  +4.92%  "}
    @{Cell = "D6"; Value = "557.59"}
    @{Cell = "E6"; Value = "  +0.32%  "}
    @{Cell = "D7"; Value = "0.999"}
    @{Cell = "E7"; Value = "  -0.15%  "}
    @{Cell = "D8"; Value = "0.583"}
    @{Cell = "E8"; Value = "  -0.75%  "}
    @{Cell = "D9"; Value = "3.313.78"}
    @{Cell = "E9"; Value = "  +0.49%  "}
    @{Cell = "D10"; Value = "0.184"}
    @{Cell = "E10"; Value = "  +0.11%  "}
    @{Cell = "D11"; Value = "0.586"}
    @{Cell = "E11"; Value = "  +0.68%  "}
    @{Cell = "D12"; Value = "47.42"}
    @{Cell = "E12"; Value = "  +0.88%  "}
    @{Cell = "E13"; Value = "  +2.82%  "}
    @{Cell = "D14"; Value = "8.70"}
    @{Cell = "E14"; Value = "  +2.46%  "}
    @{Cell = "D15"; Value = "3.850.51"}
    @{Cell = "E15"; Value = "  +0.31%  "}
    @{Cell = "D16"; Value = "608.39"}
    @{Cell = "E16"; Value = "  +1.78%  "}
    @{Cell = "D17"; Value = "66.101.23"}
    @{Cell = "E17"; Value = "  +0.89%  "}
    @{Cell = "D18"; Value = "18.02"}
    @{Cell = "E18"; Value = "  +0.17%  "}
    @{Cell = "E19"; Value = "  +1.20%  "}
    @{Cell = "D20"; Value = "3.307.36"}
    @{Cell = "E20"; Value = "  +0.12%  "}
    @{Cell = "D21"; Value = "11.11"}
    @{Cell = "E21"; Value = "  -2.35%  "}
    @{Cell = "D22"; Value = "0.909"}
    @{Cell = "E22"; Value = "  +1.43%  "}
    @{Cell = "D23"; Value = "18.46"}
    @{Cell = "E23"; Value = "  +7.79%  "}
    @{Cell = "D24"; Value = "5.08"}
    @{Cell = "E24"; Value = "  +1.15%  "}
    @{Cell = "D25"; Value = "100.48"}
    @{Cell = "E25"; Value = "  -1.81%  "}
    @{Cell = "D26"; Value = "3.97"}
    @{Cell = "E26"; Value = "  +0.30%  "}
    @{Cell = "D27"; Value = "2.78"}
    @{Cell = "E27"; Value = "  +4.59%  "}
    @{Cell = "D28"; Value = "5.94"}
    @{Cell = "E28"; Value = "  -0.90%  "}
    @{Cell = "D29"; Value = "9.60"}
    @{Cell = "E29"; Value = "  +4.06%  "}
    @{Cell = "D30"; Value = "8.73"}
    @{Cell = "E30"; Value = "  +1.26%  "}
    @{Cell = "D31"; Value = "30.41"}
    @{Cell = "E31"; Value = "  +0.04%  "}
    @{Cell = "D32"; Value = "6.76"}
    @{Cell = "E32"; Value = "  +9.20%  "}
    @{Cell = "D33"; Value = "3.89"}
    @{Cell = "E33"; Value = "  +1.93%  "}
    @{Cell = "D34"; Value = "583.20"}
    @{Cell = "E34"; Value = "  +12.45%  "}
    @{Cell = "D35"; Value = "11.09"}
    @{Cell = "E35"; Value = "  +1.09%  "}
    @{Cell = "E36"; Value = "  +1.07%  "}
    @{Cell = "D37"; Value = "3.713.87"}
    @{Cell = "E37"; Value = "  -2.17%  "}
    @{Cell = "D38"; Value = "1.00"}
    @{Cell = "E38"; Value = "  -0.02%  "}
    @{Cell = "E39"; Value = "  +0.92%  "}
    @{Cell = "D40"; Value = "34.02"}
    @{Cell = "E40"; Value = "  +6.96%  "}
    @{Cell = "B41"; Value = "Kaspa"}
    @{Cell = "C41"; Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"}
    @{Cell = "D41"; Value = "0.131"}
    @{Cell = "E41"; Value = "  +6.22%  "}
    @{Cell = "B42"; Value = "PEPE"}
    @{Cell = "C42"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"}
    @{Cell = "D42"; Value = "0.0₃0719"}
    @{Cell = "E42"; Value = "  +1.94%  "}
    @{Cell = "D43"; Value = "3.29"}
    @{Cell = "E43"; Value = "  -3.98%  "}
    @{Cell = "E44"; Value = "  +3.97%  "}
    @{Cell = "D45"; Value = "2.68"}
    @{Cell = "E45"; Value = "  +2.40%  "}
    @{Cell = "D46"; Value = "0.341"}
    @{Cell = "E46"; Value = "  +1.40%  "}
    @{Cell = "E47"; Value = "  +3.19%  "}
    @{Cell = "D48"; Value = "0.0422"}
    @{Cell = "E48"; Value = "  +3.11%  "}
    @{Cell = "E49"; Value = "  +0.59%  "}
    @{Cell = "D50"; Value = "2.60"}
    @{Cell = "E50"; Value = "  +0.66%  "}
    @{Cell = "D51"; Value = "0.998"}
    @{Cell = "E51"; Value = "  +0.04%  "}
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $change.Value
}
